# Freq_all.xlsx edit: rename Sheet1 -> data, add a MetaData sheet with
# a small "About this workbook" summary, and tidy a couple of workbook-
# level bits (active tab / calc id).

$wb = $excel.ActiveWorkbook

# --- rename the data sheet -------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "data"

# --- add the MetaData sheet right after "data" ------------------------------
$wsMeta = $wb.Worksheets.Add($null, $wsData)
$wsMeta.Name = "MetaData"

# Row 1: big title
$wsMeta.Range("A1").Value = "MetaData"
$wsMeta.Range("A1:D1").Style = "Title"
$wsMeta.Range("A1:D1").HorizontalAlignment = -4108
$wsMeta.Range("A1:D1").WrapText = $true
$wsMeta.Range("A1:D1").Merge()
$wsMeta.Rows.Item(1).RowHeight = 24

# Row 2: description of the Node column
$wsMeta.Range("A2").Value = "Node: The clusters with identifiers 0-2499 that are involved in the solution for that trial."
$wsMeta.Range("A2:D2").HorizontalAlignment = -4131
$wsMeta.Range("A2:D2").WrapText = $true
$wsMeta.Range("A2:D2").Merge()
$wsMeta.Rows.Item(2).RowHeight = 32

# Row 3: "Results" section heading
$wsMeta.Range("A3").Value = "Results"
$wsMeta.Range("A3:D3").Style = "Heading 1"
$wsMeta.Range("A3:D3").HorizontalAlignment = -4108
$wsMeta.Range("A3:D3").WrapText = $true
$b3 = $wsMeta.Range("A3:D3").Borders.Item(9)
$b3.Weight = 4
$b3.Color = 12874308
$wsMeta.Range("A3:D3").Merge()
$wsMeta.Rows.Item(3).RowHeight = 21

# Row 4: description of the f# columns
$wsMeta.Range("A4").Value = "f#: The frequency of nodes for a trial # (1,2,3,4 correspond to 1,2,2+,2-)"
$wsMeta.Range("A4:D4").HorizontalAlignment = -4131
$wsMeta.Range("A4:D4").WrapText = $true
$wsMeta.Range("A4:D4").Merge()
$wsMeta.Rows.Item(4).RowHeight = 34

# Row 5: description of the fall column
$wsMeta.Range("A5").Value = "fall: sum of frequencies of nodes over all 4 models"
$wsMeta.Range("A5:D5").HorizontalAlignment = -4131
$wsMeta.Range("A5:D5").WrapText = $true
$wsMeta.Range("A5:D5").Merge()

# Row 6: "Calculation" section heading
$wsMeta.Range("A6").Value = "Calculation"
$wsMeta.Range("A6:D6").Style = "Heading 1"
$wsMeta.Range("A6:D6").HorizontalAlignment = -4108
$wsMeta.Range("A6:D6").WrapText = $true
$b6 = $wsMeta.Range("A6:D6").Borders.Item(9)
$b6.Weight = 4
$b6.Color = 12874308
$wsMeta.Range("A6:D6").Merge()
$wsMeta.Rows.Item(6).RowHeight = 21

# Row 7: description of the r column
$wsMeta.Range("A7").Value = "r: the row of that cluster, for plotting (0,49)"
$wsMeta.Range("A7:D7").HorizontalAlignment = -4108
$b7 = $wsMeta.Range("A7:D7").Borders.Item(8)
$b7.Weight = 4
$b7.Color = 12874308
$wsMeta.Range("A7:D7").Merge()
$wsMeta.Rows.Item(7).RowHeight = 17

# Row 8: description of the c column
$wsMeta.Range("A8").Value = "c: the column of that cluster, for plotting (0-49)"
$wsMeta.Range("A8:D8").HorizontalAlignment = -4108
$wsMeta.Range("A8:D8").Merge()

$wsMeta.Range("F4").Select()

# --- workbook-level housekeeping -------------------------------------------
$wb.Windows.Item(1).ActiveSheet = $wsMeta
